$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'bjj compression shorts'
$ws.Range('A2').Value = 'compression running pants'
$ws.Range('A3').Value = 'sit pad'
$ws.Range('A4').Value = 'paintball pants padded'
$ws.Range('A5').Value = 'knee pads breathable'
$ws.Range('A6').Value = 'basketball compression tights for women'
$ws.Range('A7').Value = 'elastic band black mountain'
$ws.Range('A8').Value = 'lacrosse pads youth boys'
$ws.Range('A9').Value = 'position pad'
$ws.Range('A10').Value = 'knee length tights'
$ws.Range('A11').Value = 'capri pants men'
$ws.Range('A12').Value = 'knee pads volleyball black'
$ws.Range('A13').Value = 'softball sliding pants women'
$ws.Range('A14').Value = '6ft basketball'
$ws.Range('A15').Value = 'basketball shorts and pants'
$ws.Range('A16').Value = 'work pants for men construction knee pads'
$ws.Range('A17').Value = 'sliding shorts women softball'
$ws.Range('A18').Value = 'fight shorts wrestling'
$ws.Range('A19').Value = 'youth mesh leggings'
$ws.Range('A20').Value = 'above the knee basketball shorts'
$ws.Range('A21').Value = 'mens yoga leggings'
$ws.Range('A22').Value = 'weight lifting pants for men'
$ws.Range('A23').Value = 'cheap volleyball knee pads'
$ws.Range('A24').Value = 'compression spandex'
$ws.Range('A25').Value = 'yoga pants compression'
$ws.Range('A26').Value = 'mens above the knee shorts'
$ws.Range('A27').Value = 'mens running compression'
$ws.Range('A28').Value = 'cycling pants mens'
$ws.Range('A29').Value = 'knee sleeves basketball men'
$ws.Range('A30').Value = 'softball gear for girls'
$ws.Range('A31').Value = 'cold knee compression'
$ws.Range('A32').Value = 'youth padded compression shorts'
$ws.Range('A33').Value = 'yoga pants for men'
$ws.Range('A34').Value = 'mens spandex tights'
$ws.Range('A35').Value = 'softball protective gear'
$ws.Range('A36').Value = 'soccer sliding shorts'
$ws.Range('A37').Value = 'compression baseball shorts'
$ws.Range('A38').Value = 'long shorts for men below knee'
$ws.Range('A39').Value = 'padded leggings for cycling'
$ws.Range('A40').Value = 'padded volleyball shorts'
$ws.Range('A41').Value = 'hex squat'
$ws.Range('A42').Value = 'youth padded sliding shorts'
$ws.Range('A43').Value = 'knee sleeves bjj'
$ws.Range('A44').Value = 'football pants pads adult'
$ws.Range('A45').Value = 'work pants knee'
$ws.Range('A46').Value = 'cold compression knee'
$ws.Range('A47').Value = '5 pad football girdle'
$ws.Range('A48').Value = 'wrestling sleeve youth'
$ws.Range('A49').Value = 'compression sports pants'
$ws.Range('A50').Value = 'basketball tights for girls'
$ws.Range('A51').Value = 'water pants'
$ws.Range('A52').Value = 'spandex tights men'
$ws.Range('A53').Value = 'boys compression pants black'
$ws.Range('A54').Value = 'hockey tights'
$ws.Range('A55').Value = 'youth hockey compression pants'
$ws.Range('A56').Value = 'men leggings compression'
$ws.Range('A57').Value = 'wrestling kneepads'
$ws.Range('A58').Value = 'kneeling pad gym'
$ws.Range('A59').Value = 'guard shorts'
$ws.Range('A60').Value = 'padded compression shorts men'
$ws.Range('A61').Value = 'softball pants youth'
$ws.Range('A62').Value = 'spandex basketball shorts'
$ws.Range('A63').Value = 'compression shorts men 5 pack'
$ws.Range('A64').Value = 'shorts for men below knee'
$ws.Range('A65').Value = 'mens gym leggings'
$ws.Range('A66').Value = 'compression running leggings'
$ws.Range('A67').Value = 'black mens basketball shorts'
$ws.Range('A68').Value = 'knee pads impact'
$ws.Range('A69').Value = 'paintball pads'
$ws.Range('A70').Value = 'boys compression'
$ws.Range('A71').Value = 'mens volleyball kneepads'
$ws.Range('A72').Value = 'yoga knee pads'
$ws.Range('A73').Value = 'knee work pads'
$ws.Range('A74').Value = 'running capri'
$ws.Range('A75').Value = 'paintball pants for men'
$ws.Range('A76').Value = 'kneepad youth'
$ws.Range('A77').Value = 'polyester capri pants'
$ws.Range('A78').Value = 'man capri pants'
$ws.Range('A79').Value = 'indoor baseball'
$ws.Range('A80').Value = 'softball compression sleeve'
$ws.Range('A81').Value = 'male pads'
$ws.Range('A82').Value = 'high five girls softball pants'
$ws.Range('A83').Value = 'outdoor hockey pants'
$ws.Range('A84').Value = 'basketball floor'
$ws.Range('A85').Value = 'basketball knee sleeve black'
$ws.Range('A86').Value = 'long shorts for men below knee sports'
$ws.Range('A87').Value = 'knee pads for adults'
$ws.Range('A88').Value = 'hockey leggings'
$ws.Range('A89').Value = 'volleyball long knee pads'
$ws.Range('A90').Value = 'lacrosse shorts mens'
$ws.Range('A91').Value = 'mens tights with pouch'
$ws.Range('A92').Value = 'black short baseball pants'
$ws.Range('A93').Value = 'lightweight knee pads'
$ws.Range('A94').Value = 'mens compression pants cold'
$ws.Range('A95').Value = 'knee shorts'
$ws.Range('A96').Value = 'girls sliding pants'
$ws.Range('A97').Value = 'knee pads for work for men'
$ws.Range('A98').Value = 'youth padded leg sleeves for basketball'
$ws.Range('A99').Value = 'gym knee compression'
$ws.Range('A100').Value = 'compression football girdle'
